# Add a new "Netherlands" market sheet, cloned from the existing "Spain"
# sheet (same layout/styles), fill in its market-specific values, and
# update tab selection / active-sheet state to match.

$wb = $excel.ActiveWorkbook

# The "Spain" sheet is the template for every country tab in this workbook.
$spain = $wb.Worksheets.Item("Spain")

# Copy it, placing the new copy right after "Spain". Excel names the copy
# "Spain (2)" automatically.
$spain.Copy($null, $spain)
$netherlands = $wb.Worksheets.Item($spain.Index + 1)
$netherlands.Name = "Netherlands"

# Fill in the market-specific cells (B4 first so the shared-string table
# picks up "NGC-3144/T2034" before "Netherlands Market", matching the
# order a human typing User Story then Description would produce).
$netherlands.Range("B4").Value = "NGC-3144/T2034"
$netherlands.Range("B2").Value = "Netherlands Market"

# Restore Spain's selection to a full-sheet selection (its tab is no
# longer the active one), then make the new sheet active with B2 selected.
$spain.Cells.Select() | Out-Null
$netherlands.Activate() | Out-Null
$netherlands.Range("B2").Select() | Out-Null
